$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("http://purl.obolibrary.org/obo/OBI_0002614", "birth cohort study design", "y"),
    @("http://purl.obolibrary.org/obo/OBI_0002615", "disease specific study design", "y"),
    @("http://purl.obolibrary.org/obo/OBI_0002618", "national biomedical registry", "y"),
    @("http://purl.obolibrary.org/obo/OBI_0002617", "national registry", "y"),
    @("http://purl.obolibrary.org/obo/OBI_0002616", "genealogical record", "y")
)

$startRow = 272
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# First imported row carries wrapped text formatting on column B (matches the
# existing pattern used elsewhere in the sheet, e.g. row 156) and the taller
# row height that comes with it.
$ws.Range("B272").WrapText = $True
$ws.Rows.Item(272).RowHeight = 16

# Move the selection onto the newly imported block, like Excel leaves it
# after pasting/entering the new rows.
$ws.Range("A272:A276").Select() | Out-Null
